$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 6
$ws.Range("B14").Value = 7

[void]$ws.Range("B14").Select()
